$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap store names: row 3 (was "Bibi Cell Vieiralves") becomes "Bibi Cell Manauara",
# --- row 5 (was "Bibi Cell Manauara") becomes "Bibi Cell Vieiralves".
$ws.Range("A3").Value = "Bibi Cell Manauara"
$ws.Range("A5").Value = "Bibi Cell Vieiralves"

# --- Row 2 updates ---
$ws.Range("E2").Value = 14184.24
$ws.Range("F2").Value = 1995.8
$ws.Range("AG2").Value = 46934.31

# --- Row 3 updates ---
$ws.Range("B3").Value = 2756
$ws.Range("C3").Value = 3763
$ws.Range("D3").Value = 2753
$ws.Range("E3").Value = 2701
$ws.Range("F3").Value = 6020
$ws.Range("G3").Value = 3870.9
$ws.Range("AG3").Value = 21863.9

# --- Row 4 updates ---
$ws.Range("E4").Value = 5469
$ws.Range("F4").Value = 3040.95
$ws.Range("G4").Value = 1795.75
$ws.Range("AG4").Value = 19660.52

# --- Row 5 updates ---
$ws.Range("B5").Value = 3638
$ws.Range("C5").Value = 3280.25
$ws.Range("D5").Value = 5521.8
$ws.Range("E5").Value = 2850
$ws.Range("F5").Value = 4180
$ws.Range("AG5").Value = 19470.05

# --- Row 6 updates ---
$ws.Range("D6").Value = 22892.33
$ws.Range("E6").Value = 25204.24
$ws.Range("F6").Value = 15236.75
$ws.Range("G6").Value = 5666.65
$ws.Range("AG6").Value = 107928.78
